$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oregon's "Management actions" cell (C8) is updated to the same fuller
# description already used for Washington (D8), noting that the
# evisceration order only applied in 2015, and picks up the matching
# wrap-text style.
$ws.Range("C8").Value = "Area closure, evisceration order" + [char]10 + "(area closure only in 2015)"
$ws.Range("C8").WrapText = $true

# Washington's cell (D8) now shares that same string (no visible change).
$ws.Range("D8").Value = "Area closure, evisceration order" + [char]10 + "(area closure only in 2015)"
$ws.Range("D8").WrapText = $true

# Update the saved selection/active cell.
[void]$ws.Range("C13").Select()
